# Adds a new "Info" worksheet (placed after "ColumnsNumberParameters") that
# documents the purpose of the workbook, with a couple of words/phrases
# emphasised via italic / italic+underline rich-text runs - mirroring the
# target commit's new Info tab.

$wb = $excel.ActiveWorkbook
$wsMain = $wb.Worksheets.Item("ColumnsNumberParameters")

# New sheet is inserted right after the main sheet and becomes the active one
# (Excel automatically flips bookViews/workbookView@activeTab + the new
# sheet's tabSelected when a freshly added sheet becomes active).
$wsInfo = $wb.Worksheets.Add($null, $wsMain)
$wsInfo.Name = "Info"

# ---------------------------------------------------------------------
# A1: "... in the reg_estimates files." - "reg_estimates" set to italic
# ---------------------------------------------------------------------
$text1 = "This Excel file is used to define the column numbers required for the corresponding processes in the reg_estimates files."
$wsInfo.Range("A1").Value = $text1

$kw1 = "reg_estimates"
$start1 = $text1.IndexOf($kw1) + 1
$len1 = $kw1.Length

$wsInfo.Range("A1").Characters($start1, $len1).Font.Italic = $true

$tailStart1 = $start1 + $len1
$tailLen1 = $text1.Length - $tailStart1 + 1
$tail1 = $wsInfo.Range("A1").Characters($tailStart1, $tailLen1)
$tail1.Font.Italic = $false

# ---------------------------------------------------------------------
# A2: "... configured separately for each country." - "separately" and
# "each country" set to italic + underline
# ---------------------------------------------------------------------
$text2 = "Since the column numbers differ by country, they must be configured separately for each country."
$wsInfo.Range("A2").Value = $text2

$kw2a = "separately"
$start2a = $text2.IndexOf($kw2a) + 1
$len2a = $kw2a.Length
$run2a = $wsInfo.Range("A2").Characters($start2a, $len2a)
$run2a.Font.Italic = $true
$run2a.Font.Underline = $true

$kw2b = "each country"
$start2b = $text2.IndexOf($kw2b) + 1
$len2b = $kw2b.Length

$midStart2 = $start2a + $len2a
$midLen2 = $start2b - $midStart2
$mid2 = $wsInfo.Range("A2").Characters($midStart2, $midLen2)
$mid2.Font.Italic = $false
$mid2.Font.Underline = $false

$run2b = $wsInfo.Range("A2").Characters($start2b, $len2b)
$run2b.Font.Italic = $true
$run2b.Font.Underline = $true

$tailStart2 = $start2b + $len2b
$tailLen2 = $text2.Length - $tailStart2 + 1
$tail2 = $wsInfo.Range("A2").Characters($tailStart2, $tailLen2)
$tail2.Font.Italic = $false
$tail2.Font.Underline = $false

# Match the saved selection on the new Info sheet (A1:A2)
[void]$wsInfo.Range("A1:A2").Select()

Write-Output "Info sheet added"
